$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 13 ("Docentes responsáveis:" value row).
#    This shifts the previous rows 13-23 down to 14-24, carrying their
#    row-height / customHeight formatting with them (already matches target).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# ---------------------------------------------------------------------------
# 2) Row 10 ("Objetivos:") - fix the stray leftover value in B/C and replace
#    it with the real course-objectives text.
# ---------------------------------------------------------------------------
$objetivos = @"
Levar ao aluno uma visão relativamente aprofundada sobre a ciência dos polímeros. Apresentar os conceitos fundamentais, os mecanismos envolvidos nas sínteses dos polímeros, os diferentes processos de polimerização e finalmente as propriedades mais marcantes dos materiais obtidos. Mostrar para o aluno a importância do conhecimento destes materiais na sua carreira profissional.
"@
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------------
# 3) New row 13 ("Docentes responsáveis:" value) - no label in column A,
#    value in B/C is the professor's name. The freshly-inserted row has no
#    column formatting yet, so copy B/C formats down from row 14 (which
#    already carries the correct "value column" styles) before writing.
# ---------------------------------------------------------------------------
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$docentes = "5840772 - Amilton Martins dos Santos"
$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes

# ---------------------------------------------------------------------------
# 4) Row 14 ("Programa resumido:") - replace placeholder "Semestral" text
#    with the real short-syllabus summary.
# ---------------------------------------------------------------------------
$resumido = @"
Introdução a polímeros; Mecanismos de polimerização; Técnicas de polimerização; Processos de polimerização; Caracterização de polímeros; Propriedades de polímeros.
"@
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# ---------------------------------------------------------------------------
# 5) Row 16 ("Programa:") - fix stray leftover value, add the full syllabus.
# ---------------------------------------------------------------------------
$programa = @"
Introdução a polímeros: História, Conceitos fundamentais, Classificação dos Polímeros, Nomenclatura de polímeros. Mecanismos de polimerização: Definições de poliadição e policondensação, policondensação (poliésteres, poliamidas, policarbonatos, poliuretanos), poliadição (polimerização via radical livre). Técnicas de polimerização (Massa, solução, suspensão, emulsão e miniemulsão). Processos de polimerização (Batelada, batelada alimentada/semi-contínuo, processo contínuo, processo shot). Noções sobre a caracterização de polímeros (GPC/SEC, DSC e TGA). Definição das propriedades mais importantes dos polímeros (Tg e TM, outras propriedades de engenharia).
"@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# ---------------------------------------------------------------------------
# 6) Row 19 ("Método:") - fix stray leftover value, add the evaluation method.
# ---------------------------------------------------------------------------
$metodo = "2 Provas escritas + Trabalho de conclusão de curso."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# ---------------------------------------------------------------------------
# 7) Row 20 ("Critério:") - now holds the grading-formula text.
# ---------------------------------------------------------------------------
$criterio = @"
A nota final (NF) será calculada de seguinte maneira: NF = (P1+P2)/2
O trabalho poderá valer até 2 pontos, que serão somados nas notas da P1 ou da P2.
"@
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# ---------------------------------------------------------------------------
# 8) Row 21 ("Norma de recuperação:") - now holds the recovery-exam text.
# ---------------------------------------------------------------------------
$norma = @"
A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula MR = (NF+PR)/2.
"@
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# ---------------------------------------------------------------------------
# 9) Row 22 ("Bibliografia:") - replace with the real bibliography list.
# ---------------------------------------------------------------------------
$biblio = @"
MANO E. B. Introdução a Polímeros. Editora Edgard Blücher Ltda, 1a Ed., São Paulo, 1988; MANO E. B. Polímeros como Materiais de Engenharia. Editora Edgard Blücher Ltda, 1a Ed., São Paulo, 1991
CANNEVALORO S. V. Ciência dos Polímeros. Editora Artliber  Ltda, 1a Ed., São Paulo, 2004
COUTINHO F. M. B.; OLIVEIRA C. M. F. Reações de Polimerização em Cadeia. Editora Interciência Ltda, 1ª Ed., Rio de Janeiro, 2006
BILMEYER Jr., F. W. Textbook of Polymer Science. John Wiley & Sons, 3rd Ed., New York, 1984
ODIAN G. Principles of Polymerization, John Wiley & Sons, 3rd Ed., New York, 1991
RODRIGUEZ, FERDINAND. Princípios de Sistemas de Polímeros, Editorial El Manual Moderno S.A., 1st Ed., México, D.F., 1984
VAN KREVELEN, D.W., HOFTYZER, P. J. Properties of polymers: correlation with chemical structure. Elsevier, 1st Ed., Amsterdam, 1972.
"@
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

# ---------------------------------------------------------------------------
# 10) Column layout clean-up: column A's width definition should cover only
#     column A (it previously also overlapped column B's range).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
